$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.034771960228681564
$ws.Range("C2").Value = 0.015054875053465366
$ws.Range("D2").Value = 0.01084665022790432
$ws.Range("E2").Value = 0.007411157246679068
$ws.Range("F2").Value = [double]"7.001341873547062E-5"
$ws.Range("J2").Value = 0.12762977182865143
$ws.Range("K2").Value = 1.459633708000183
